# Auto-generated script to update Leve profit calculation cells
# across multiple worksheets, reflecting refreshed market board prices.
$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC (59 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2916.6667
$ws.Range("J17").Value = 2916.6667
$ws.Range("L17").Value = 8750.000100000001
$ws.Range("N17").Value = -9086.000100000001
$ws.Range("H18").Value = 1248.3334
$ws.Range("I18").Value = 1248.3334
$ws.Range("K18").Value = 1248.3334
$ws.Range("M18").Value = -964.3334
$ws.Range("H33").Value = 14019.277
$ws.Range("I33").Value = 15728.625
$ws.Range("K33").Value = 15728.625
$ws.Range("M33").Value = -15499.625
$ws.Range("H40").Value = 1486
$ws.Range("I40").Value = 1486
$ws.Range("K40").Value = 1486
$ws.Range("M40").Value = -1311
$ws.Range("H70").Value = 4356
$ws.Range("I70").Value = 4374.75
$ws.Range("J70").Value = 4331
$ws.Range("K70").Value = 13124.25
$ws.Range("L70").Value = 12993
$ws.Range("M70").Value = -12854.25
$ws.Range("N70").Value = -13533
$ws.Range("H73").Value = 4356
$ws.Range("I73").Value = 4374.75
$ws.Range("J73").Value = 4331
$ws.Range("K73").Value = 13124.25
$ws.Range("L73").Value = 12993
$ws.Range("M73").Value = -12188.25
$ws.Range("N73").Value = -14865
$ws.Range("H100").Value = 2982.6365
$ws.Range("I100").Value = 2162.25
$ws.Range("J100").Value = 5170.3335
$ws.Range("K100").Value = 2162.25
$ws.Range("L100").Value = 5170.3335
$ws.Range("M100").Value = -1621.25
$ws.Range("N100").Value = -6252.3335
$ws.Range("H112").Value = 3679.9412
$ws.Range("J112").Value = 3681.8125
$ws.Range("L112").Value = 11045.4375
$ws.Range("N112").Value = -13261.4375
$ws.Range("H135").Value = 1856.3529
$ws.Range("I135").Value = 504.91666
$ws.Range("J135").Value = 5099.8
$ws.Range("K135").Value = 4544.24994
$ws.Range("L135").Value = 45898.2
$ws.Range("M135").Value = -2009.24994
$ws.Range("N135").Value = -50968.2
$ws.Range("H137").Value = 23261142
$ws.Range("I137").Value = 45455760
$ws.Range("K137").Value = 136367280
$ws.Range("M137").Value = -136364730
$ws.Range("H141").Value = 2987.6191
$ws.Range("I141").Value = 2743.625
$ws.Range("J141").Value = 3137.7693
$ws.Range("K141").Value = 8230.875
$ws.Range("L141").Value = 9413.3079
$ws.Range("M141").Value = -3050.875
$ws.Range("N141").Value = -19773.3079

# --- Worksheet: ARM (40 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1551.65
$ws.Range("I45").Value = 1431.2858
$ws.Range("J45").Value = 1832.5
$ws.Range("K45").Value = 1431.2858
$ws.Range("L45").Value = 1832.5
$ws.Range("M45").Value = -1054.2858
$ws.Range("N45").Value = -2586.5
$ws.Range("H61").Value = 3338266
$ws.Range("I61").Value = 5043.727
$ws.Range("J61").Value = 12504627
$ws.Range("K61").Value = 5043.727
$ws.Range("L61").Value = 12504627
$ws.Range("M61").Value = -4831.727
$ws.Range("N61").Value = -12505051
$ws.Range("H74").Value = 900836.2
$ws.Range("I74").Value = 1070132.6
$ws.Range("K74").Value = 1070132.6
$ws.Range("M74").Value = -1069258.6
$ws.Range("H77").Value = 900836.2
$ws.Range("I77").Value = 1070132.6
$ws.Range("K77").Value = 5350663
$ws.Range("M77").Value = -5346295
$ws.Range("H110").Value = 761.4666999999999
$ws.Range("I110").Value = 569.96
$ws.Range("K110").Value = 569.96
$ws.Range("M110").Value = 1475.04
$ws.Range("H132").Value = 326871.38
$ws.Range("I132").Value = 359157.16
$ws.Range("J132").Value = 4013.5715
$ws.Range("K132").Value = 1077471.48
$ws.Range("L132").Value = 12040.7145
$ws.Range("M132").Value = -1074941.48
$ws.Range("N132").Value = -17100.7145
$ws.Range("H136").Value = 3338266
$ws.Range("I136").Value = 5043.727
$ws.Range("J136").Value = 12504627
$ws.Range("K136").Value = 15131.181
$ws.Range("L136").Value = 37513881
$ws.Range("M136").Value = -12581.181
$ws.Range("N136").Value = -37518981

# --- Worksheet: CRP (54 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1636408.6
$ws.Range("I31").Value = 2059889.1
$ws.Range("J31").Value = 2983.8572
$ws.Range("K31").Value = 2059889.1
$ws.Range("L31").Value = 2983.8572
$ws.Range("M31").Value = -2059594.1
$ws.Range("N31").Value = -3573.8572
$ws.Range("H34").Value = 1636408.6
$ws.Range("I34").Value = 2059889.1
$ws.Range("J34").Value = 2983.8572
$ws.Range("K34").Value = 2059889.1
$ws.Range("L34").Value = 2983.8572
$ws.Range("M34").Value = -2059687.1
$ws.Range("N34").Value = -3387.8572
$ws.Range("H99").Value = 31971.5
$ws.Range("I99").Value = 38499.816
$ws.Range("J99").Value = 17609.2
$ws.Range("K99").Value = 38499.816
$ws.Range("L99").Value = 17609.2
$ws.Range("M99").Value = -37001.816
$ws.Range("N99").Value = -20605.2
$ws.Range("H107").Value = 890
$ws.Range("I107").Value = 927.6667
$ws.Range("K107").Value = 927.6667
$ws.Range("M107").Value = 992.3333
$ws.Range("H122").Value = 7914.615
$ws.Range("J122").Value = 43626.25
$ws.Range("L122").Value = 130878.75
$ws.Range("N122").Value = -135778.75
$ws.Range("H126").Value = 31971.5
$ws.Range("I126").Value = 38499.816
$ws.Range("J126").Value = 17609.2
$ws.Range("K126").Value = 115499.448
$ws.Range("L126").Value = 52827.60000000001
$ws.Range("M126").Value = -113029.448
$ws.Range("N126").Value = -57767.60000000001
$ws.Range("H132").Value = 2804.6667
$ws.Range("I132").Value = 3016.353
$ws.Range("J132").Value = 2444.8
$ws.Range("K132").Value = 9049.059000000001
$ws.Range("L132").Value = 7334.400000000001
$ws.Range("M132").Value = -6519.059000000001
$ws.Range("N132").Value = -12394.4
$ws.Range("H134").Value = 1817.72
$ws.Range("I134").Value = 1714.9131
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 5144.7393
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -2609.7393
$ws.Range("N134").Value = -14070
$ws.Range("H141").Value = 217438.53
$ws.Range("J141").Value = 225827.08
$ws.Range("L141").Value = 225827.08
$ws.Range("N141").Value = -236187.08

# --- Worksheet: CUL (4 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 332.5
$ws.Range("I107").Value = 449.66666
$ws.Range("K107").Value = 1348.99998
$ws.Range("M107").Value = 571.0000199999999

# --- Worksheet: GSM (20 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 24864.25
$ws.Range("J20").Value = 23985.666
$ws.Range("L20").Value = 23985.666
$ws.Range("N20").Value = -24475.666
$ws.Range("H24").Value = 19575
$ws.Range("J24").Value = 21005
$ws.Range("L24").Value = 21005
$ws.Range("N24").Value = -21351
$ws.Range("H123").Value = 44325.11
$ws.Range("J123").Value = 44325.11
$ws.Range("L123").Value = 44325.11
$ws.Range("N123").Value = -49225.11
$ws.Range("H132").Value = 12068.189
$ws.Range("I132").Value = 10110.061
$ws.Range("K132").Value = 30330.183
$ws.Range("M132").Value = -27800.183
$ws.Range("H136").Value = 82701.37
$ws.Range("J136").Value = 82701.37
$ws.Range("L136").Value = 248104.11
$ws.Range("N136").Value = -253204.11

# --- Worksheet: LTW (7 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2909286.8
$ws.Range("I132").Value = 3934635
$ws.Range("J132").Value = 4133.3335
$ws.Range("K132").Value = 11803905
$ws.Range("L132").Value = 12400.0005
$ws.Range("M132").Value = -11801375
$ws.Range("N132").Value = -17460.0005

# --- Worksheet: WVR (20 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 24709.166
$ws.Range("I31").Value = 24000
$ws.Range("J31").Value = 24851
$ws.Range("K31").Value = 24000
$ws.Range("L31").Value = 24851
$ws.Range("M31").Value = -23652
$ws.Range("N31").Value = -25547
$ws.Range("H51").Value = 15821.154
$ws.Range("I51").Value = 12144.4
$ws.Range("K51").Value = 12144.4
$ws.Range("M51").Value = -11634.4
$ws.Range("H52").Value = 23535.25
$ws.Range("H107").Value = 2469.4443
$ws.Range("I107").Value = 2743.75
$ws.Range("K107").Value = 8231.25
$ws.Range("M107").Value = -6311.25
$ws.Range("H136").Value = 8117741
$ws.Range("I136").Value = 2071753.6
$ws.Range("K136").Value = 6215260.800000001
$ws.Range("M136").Value = -6212710.800000001

Write-Host "Updated 204 cells across 7 worksheets"